$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "domingo" cell in row 2 with the new string
$ws.Range("A2").Value = "putAKeepAreYou"

# Replace the "Total da carga (Kg)" label with "nada para nada"
$ws.Range("A7").Value = "nada para nada"

# Update the "domingo" cell in row 3 with the new string
$ws.Range("A3").Value = "´pourrra"

# Move the active selection to A12, matching the saved view state
$ws.Range("A12").Select()
